# Update the LR-pairs sheet with newly recomputed TPM-based values.
# The data now only keeps the "ECs" target-cluster rows (the "MuSCs" rows
# are dropped entirely), leaving 3 data rows instead of 6, and several
# numeric columns are recalculated with the new TPM numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that belonged to the "MuSCs" target cluster (old rows 3, 5
# and 7) plus collapse the remaining "ECs" rows into a contiguous 2..4
# block. Deleting whole rows shifts everything below upward automatically.

# Old row 3 (FAPs / Ccl12 / Ccr4 / MuSCs) -> delete
$ws.Rows.Item(3).Delete()

# After that delete, the former row 5 (Inflammatory-Mac / MuSCs) is now row 4
$ws.Rows.Item(4).Delete()

# After that delete, the former row 7 (Resolving-Mac / MuSCs) is now row 5
$ws.Rows.Item(5).Delete()

# Now rows 2..4 hold the FAPs/Inflammatory-Mac/Resolving-Mac -> ECs triples.
# Update their numeric columns with the freshly recomputed TPM values.

# Row 2: FAPs -> Ccl12 -> Ccr4 -> ECs
$ws.Range("I2").Value = 0.005723000769734084
$ws.Range("J2").Value = 0.005723000769734084
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("S2").Value = 0.005723000769734084
$ws.Range("T2").Value = 0.005723000769734084

# Row 3: Inflammatory-Mac -> Ccl12 -> Ccr4 -> ECs
$ws.Range("A3").Value = "Inflammatory-Mac"
$ws.Range("D3").Value = "ECs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 30.199365
$ws.Range("H3").Value = 90.598095
$ws.Range("I3").Value = 0.6125398923302606
$ws.Range("J3").Value = 0.6125398923302606
$ws.Range("M3").Value = 0.08819666666666666
$ws.Range("N3").Value = 0.26459
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 2.66348332845
$ws.Range("R3").Value = 23.97134995605
$ws.Range("S3").Value = 0.6125398923302606
$ws.Range("T3").Value = 0.6125398923302606

# Row 4: Resolving-Mac -> Ccl12 -> Ccr4 -> ECs
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 18.820355
$ws.Range("H4").Value = 56.461065
$ws.Range("I4").Value = 0.3817371069000054
$ws.Range("J4").Value = 0.3817371069000054
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 1.659892576483333
$ws.Range("R4").Value = 14.93903318835
$ws.Range("S4").Value = 0.3817371069000054
$ws.Range("T4").Value = 0.3817371069000054
